$wb = $excel.ActiveWorkbook
$ws7 = $wb.Worksheets.Item("optimization_parameters")

# 1. Rename the "Model" label to "production_function"
$ws7.Range("A8").Value = "production_function"

# 2. Insert a new "L_curve" row after the production_function row (shifts rows 9+ down by one)
$ws7.Range("A9:E9").Insert(-4121)   # xlShiftDown
$ws7.Range("A9").Value = "L_curve"
$ws7.Range("B9").Value = 0
$ws7.Range("B9").NumberFormat = $ws7.Range("B2").NumberFormat

# 3. Delete the old "Deletion" row (now shifted down to row 17)
$ws7.Range("A17").EntireRow.Delete()

# 4. Clear the extra header cells C1:F1
$ws7.Range("C1:F1").ClearContents()

# 5. Make optimization_parameters the active/selected sheet & set its selection
$ws7.Activate()
$ws7.Range("C1:F1").Select()
